# Add a new "Croatia" market test-data sheet, cloned from the existing
# "Turkey" sheet (same layout/styles), then localize its content.

$wb = $excel.ActiveWorkbook

$turkey = $wb.Worksheets.Item("Turkey")

# Duplicate "Turkey" and place the copy right after it.
$turkey.Copy($null, $turkey)
$croatia = $wb.Worksheets.Item("Turkey (2)")
$croatia.Name = "Croatia"

# Localize the two market-specific cells on the new sheet.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2415"

# The shorter replacement text re-wraps, so the row heights settle
# differently than on the Turkey sheet it was cloned from.
$croatia.Rows.Item(3).AutoFit()
$croatia.Rows.Item(4).RowHeight = 15.6
$croatia.Rows.Item(5).AutoFit()

# Turkey is no longer the active tab/selection; Croatia (the new last
# sheet) becomes the active one.
$turkey.Range("L20").Select()
$croatia.Activate()
$croatia.Range("J19").Select()

$wb.Save()
